$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity row: copy formatting from the row above (row 18) so that
# row 19 gets the same data-row styling (bordered, centered, date format).
$ws.Range("B18:E18").Copy()
$ws.Range("B19:E19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B19").Value = 12
$ws.Range("C19").Value = "Ajustes finales página web"
$ws.Range("D19").Value = "Soraya Soto"
$ws.Range("E19").Value = (Get-Date -Year 2021 -Month 4 -Day 12 -Hour 0 -Minute 0 -Second 0)

# D21 becomes centered like B21.
$ws.Range("D21").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D21").VerticalAlignment = -4108     # xlCenter

$ws.Range("D21").Select()
